$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: force text storage (NumberFormat "@") then reset cell
# style back to Normal so no stray style index is left on the cell,
# matching the original inline-string (unstyled) cells.

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '30.030.02'
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).NumberFormat = "@"
$ws.Cells.Item(2, 5).Value = '  -0.25%  '
$ws.Cells.Item(2, 5).Style = "Normal"

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.871.19'
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = '  -2.47%  '
$ws.Cells.Item(3, 5).Style = "Normal"

# Row 4
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.002'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).NumberFormat = "@"
$ws.Cells.Item(4, 5).Value = '  +0.08%  '
$ws.Cells.Item(4, 5).Style = "Normal"

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '319.14'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = '  -3.27%  '
$ws.Cells.Item(5, 5).Style = "Normal"

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '1.000'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).NumberFormat = "@"
$ws.Cells.Item(6, 5).Value = '  -0.03%  '
$ws.Cells.Item(6, 5).Style = "Normal"

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5056'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).NumberFormat = "@"
$ws.Cells.Item(7, 5).Value = '  -3.17%  '
$ws.Cells.Item(7, 5).Style = "Normal"

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3963'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).NumberFormat = "@"
$ws.Cells.Item(8, 5).Value = '  -3.11%  '
$ws.Cells.Item(8, 5).Style = "Normal"

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.08207'
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).NumberFormat = "@"
$ws.Cells.Item(9, 5).Value = '  -3.52%  '
$ws.Cells.Item(9, 5).Style = "Normal"

# Row 10
$ws.Cells.Item(10, 5).NumberFormat = "@"
$ws.Cells.Item(10, 5).Value = '  -2.10%  '
$ws.Cells.Item(10, 5).Style = "Normal"

# Row 11
$ws.Cells.Item(11, 5).NumberFormat = "@"
$ws.Cells.Item(11, 5).Value = '  -3.10%  '
$ws.Cells.Item(11, 5).Style = "Normal"

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '23.41'
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).NumberFormat = "@"
$ws.Cells.Item(12, 5).Value = '  +4.09%  '
$ws.Cells.Item(12, 5).Style = "Normal"

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.864.32'
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).NumberFormat = "@"
$ws.Cells.Item(13, 5).Value = '  -3.70%  '
$ws.Cells.Item(13, 5).Style = "Normal"

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '6.280'
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).NumberFormat = "@"
$ws.Cells.Item(14, 5).Value = '  -2.23%  '
$ws.Cells.Item(14, 5).Style = "Normal"

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.193'
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).NumberFormat = "@"
$ws.Cells.Item(15, 5).Value = '  -3.15%  '
$ws.Cells.Item(15, 5).Style = "Normal"

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '1.002'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).NumberFormat = "@"
$ws.Cells.Item(16, 5).Value = '  +0.09%  '
$ws.Cells.Item(16, 5).Style = "Normal"

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '91.81'
$ws.Cells.Item(17, 4).Style = "Normal"
$ws.Cells.Item(17, 5).NumberFormat = "@"
$ws.Cells.Item(17, 5).Value = '  -4.03%  '
$ws.Cells.Item(17, 5).Style = "Normal"

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.00001085'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = '  -2.48%  '
$ws.Cells.Item(18, 5).Style = "Normal"

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.06443'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = '  -4.27%  '
$ws.Cells.Item(19, 5).Style = "Normal"

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '18.09'
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = '  -1.23%  '
$ws.Cells.Item(20, 5).Style = "Normal"

# Row 21
$ws.Cells.Item(21, 5).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = '  -0.02%  '
$ws.Cells.Item(21, 5).Style = "Normal"

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '30.027.98'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).NumberFormat = "@"
$ws.Cells.Item(22, 5).Value = '  -0.26%  '
$ws.Cells.Item(22, 5).Style = "Normal"

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.843'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).NumberFormat = "@"
$ws.Cells.Item(23, 5).Value = '  -2.78%  '
$ws.Cells.Item(23, 5).Style = "Normal"

# Row 24
$ws.Cells.Item(24, 5).NumberFormat = "@"
$ws.Cells.Item(24, 5).Value = '  -1.82%  '
$ws.Cells.Item(24, 5).Style = "Normal"

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.168'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).NumberFormat = "@"
$ws.Cells.Item(25, 5).Value = '  -2.51%  '
$ws.Cells.Item(25, 5).Style = "Normal"

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.087.82'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).NumberFormat = "@"
$ws.Cells.Item(26, 5).Value = '  -3.08%  '
$ws.Cells.Item(26, 5).Style = "Normal"

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '21.27'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).NumberFormat = "@"
$ws.Cells.Item(27, 5).Value = '  +0.86%  '
$ws.Cells.Item(27, 5).Style = "Normal"

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '161.04'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).NumberFormat = "@"
$ws.Cells.Item(28, 5).Value = '  +0.30%  '
$ws.Cells.Item(28, 5).Style = "Normal"

# Row 29
$ws.Cells.Item(29, 5).NumberFormat = "@"
$ws.Cells.Item(29, 5).Value = '  -9.50%  '
$ws.Cells.Item(29, 5).Style = "Normal"

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '127.32'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).NumberFormat = "@"
$ws.Cells.Item(30, 5).Value = '  -1.49%  '
$ws.Cells.Item(30, 5).Style = "Normal"

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.071'
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).NumberFormat = "@"
$ws.Cells.Item(31, 5).Value = '  -0.63%  '
$ws.Cells.Item(31, 5).Style = "Normal"

# Row 32
$ws.Cells.Item(32, 5).NumberFormat = "@"
$ws.Cells.Item(32, 5).Value = '  -2.01%  '
$ws.Cells.Item(32, 5).Style = "Normal"

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '5.940'
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).NumberFormat = "@"
$ws.Cells.Item(33, 5).Value = '  -2.49%  '
$ws.Cells.Item(33, 5).Style = "Normal"

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.687'
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).NumberFormat = "@"
$ws.Cells.Item(34, 5).Value = '  +1.56%  '
$ws.Cells.Item(34, 5).Style = "Normal"

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.02437'
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).NumberFormat = "@"
$ws.Cells.Item(35, 5).Value = '  -2.36%  '
$ws.Cells.Item(35, 5).Style = "Normal"

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.228'
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).NumberFormat = "@"
$ws.Cells.Item(36, 5).Value = '  +0.73%  '
$ws.Cells.Item(36, 5).Style = "Normal"

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.06354'
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).NumberFormat = "@"
$ws.Cells.Item(37, 5).Value = '  -4.03%  '
$ws.Cells.Item(37, 5).Style = "Normal"

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.2141'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).NumberFormat = "@"
$ws.Cells.Item(38, 5).Value = '  -3.11%  '
$ws.Cells.Item(38, 5).Style = "Normal"

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.172'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).NumberFormat = "@"
$ws.Cells.Item(39, 5).Value = '  -4.91%  '
$ws.Cells.Item(39, 5).Style = "Normal"

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '8.480'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).NumberFormat = "@"
$ws.Cells.Item(40, 5).Value = '  -4.81%  '
$ws.Cells.Item(40, 5).Style = "Normal"

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.6303'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).NumberFormat = "@"
$ws.Cells.Item(41, 5).Value = '  -3.45%  '
$ws.Cells.Item(41, 5).Style = "Normal"

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '1.216'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).NumberFormat = "@"
$ws.Cells.Item(42, 5).Value = '  -2.57%  '
$ws.Cells.Item(42, 5).Style = "Normal"

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '11.29'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).NumberFormat = "@"
$ws.Cells.Item(43, 5).Value = '  -3.09%  '
$ws.Cells.Item(43, 5).Style = "Normal"

# Row 44
$ws.Cells.Item(44, 2).NumberFormat = "@"
$ws.Cells.Item(44, 2).Value = 'EnergySwap'
$ws.Cells.Item(44, 2).Style = "Normal"
$ws.Cells.Item(44, 3).NumberFormat = "@"
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(44, 3).Style = "Normal"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '13.05'
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).NumberFormat = "@"
$ws.Cells.Item(44, 5).Value = '  -1.70%  '
$ws.Cells.Item(44, 5).Style = "Normal"

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.5916'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).NumberFormat = "@"
$ws.Cells.Item(45, 5).Value = '  -4.09%  '
$ws.Cells.Item(45, 5).Style = "Normal"

# Row 46
$ws.Cells.Item(46, 2).NumberFormat = "@"
$ws.Cells.Item(46, 2).Value = 'NEARProtocol'
$ws.Cells.Item(46, 2).Style = "Normal"
$ws.Cells.Item(46, 3).NumberFormat = "@"
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(46, 3).Style = "Normal"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '2.093'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).NumberFormat = "@"
$ws.Cells.Item(46, 5).Value = '  +0.56%  '
$ws.Cells.Item(46, 5).Style = "Normal"

# Row 47
$ws.Cells.Item(47, 2).NumberFormat = "@"
$ws.Cells.Item(47, 2).Value = 'PancakeSwap'
$ws.Cells.Item(47, 2).Style = "Normal"
$ws.Cells.Item(47, 3).NumberFormat = "@"
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(47, 3).Style = "Normal"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '3.622'
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).NumberFormat = "@"
$ws.Cells.Item(47, 5).Value = '  -3.93%  '
$ws.Cells.Item(47, 5).Style = "Normal"

# Row 48
$ws.Cells.Item(48, 2).NumberFormat = "@"
$ws.Cells.Item(48, 2).Value = 'Quant'
$ws.Cells.Item(48, 2).Style = "Normal"
$ws.Cells.Item(48, 3).NumberFormat = "@"
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(48, 3).Style = "Normal"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '122.60'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).NumberFormat = "@"
$ws.Cells.Item(48, 5).Value = '  -1.63%  '
$ws.Cells.Item(48, 5).Style = "Normal"

# Row 49
$ws.Cells.Item(49, 2).NumberFormat = "@"
$ws.Cells.Item(49, 2).Value = 'EOS'
$ws.Cells.Item(49, 2).Style = "Normal"
$ws.Cells.Item(49, 3).NumberFormat = "@"
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Cells.Item(49, 3).Style = "Normal"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.205'
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = '  -3.41%  '
$ws.Cells.Item(49, 5).Style = "Normal"

# Row 50
$ws.Cells.Item(50, 2).NumberFormat = "@"
$ws.Cells.Item(50, 2).Value = 'Aave'
$ws.Cells.Item(50, 2).Style = "Normal"
$ws.Cells.Item(50, 3).NumberFormat = "@"
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(50, 3).Style = "Normal"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '77.39'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).NumberFormat = "@"
$ws.Cells.Item(50, 5).Value = '  -3.08%  '
$ws.Cells.Item(50, 5).Style = "Normal"

# Row 51
$ws.Cells.Item(51, 2).NumberFormat = "@"
$ws.Cells.Item(51, 2).Value = 'WEMIXTOKEN'
$ws.Cells.Item(51, 2).Style = "Normal"
$ws.Cells.Item(51, 3).NumberFormat = "@"
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(51, 3).Style = "Normal"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '1.112'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).NumberFormat = "@"
$ws.Cells.Item(51, 5).Value = '  -4.39%  '
$ws.Cells.Item(51, 5).Style = "Normal"
